# Apply cryptos list update (Tue Aug 13 04:45:27 UTC 2024, GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All source values are plain text (prices/percent strings with thousands-dot
# separators, leading zeros, etc.) - force text format so Excel does not
# auto-coerce them into numbers/dates and mangle formatting (e.g. "0.0230").
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.016.07"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.634.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.72%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.32"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.75"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.657.87"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.24"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.58%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.337"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.101.17"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.951.97"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.91"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.651.96"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "347.97"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.84%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.19"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.93%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.77"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.750.66"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.77%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.994"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.161"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0805"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.38%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.14"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.46%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.27"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +7.35%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.97"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.57"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.42%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.60"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "SuiNetwork"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.974"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.88%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.01"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.53%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.14"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.74"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.04%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.844"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.56%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.70"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.83%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.41"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "278.22"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.612"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.994"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0983"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.62"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.13%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0529"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.30"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0230"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.989.41"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.35%  "
